$d = $word.ActiveDocument

$pairs = @(
    @("148÷9=", "998÷3="),
    @("418÷3=", "446÷8="),
    @("129÷8=", "508÷8="),
    @("527÷2=", "672÷6="),
    @("274÷2=", "409÷7="),
    @("678÷8=", "820÷6="),
    @("998÷4=", "218÷5="),
    @("343÷6=", "390÷4="),
    @("581÷8=", "584÷2="),
    @("959÷5=", "440÷9="),
    @("350÷6=", "923÷5="),
    @("558÷7=", "903÷5="),
    @("933÷8=", "981÷3="),
    @("188÷2=", "593÷9="),
    @("302÷4=", "825÷4="),
    @("210÷8=", "147÷9="),
    @("636÷7=", "950÷6="),
    @("907÷6=", "750÷8="),
    @("787÷3=", "797÷2="),
    @("412÷8=", "272÷9="),
    @("791÷5=", "688÷9="),
    @("600÷7=", "317÷6="),
    @("868÷3=", "917÷6="),
    @("608÷7=", "461÷7="),
    @("538÷5=", "234÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
